# Update the crypto price/volume table (columns D "Price" and E "Volume(1h)")
# for rows 2-51 to the newly scraped values.
#
# All D/E cells in the source sheet are stored as *text* (t="inlineStr"),
# e.g. "0.7120" or "  -0.85%  ", even though some look numeric. Excel's COM
# layer auto-converts plain numeric-looking strings (no "%", no double dots)
# into real numbers when assigned directly, which would lose formatting such
# as trailing zeros. To keep those cells as text we prefix the value with a
# leading apostrophe (the normal Excel "force text" marker) and then clear the
# resulting cell formatting so no stray number-format style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.227.53'
$ws.Range("D3").Value = '1.862.48'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'0.7120"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").Value = "'240.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.3082"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = "'0.07693"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").Value = "'24.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = "'0.08374"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("D12").Value = '1.856.23'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").Value = "'5.200"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("D14").Value = "'0.7119"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").Value = "'91.29"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '29.228.68'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = "'5.948"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = "'242.35"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.000007833"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = '2.121.10'
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").Value = "'13.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = "'7.833"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'0.1592"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = "'163.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Value = "'8.887"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = "'18.47"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").Value = "'1.341"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").Value = "'1.497"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").Value = "'4.404"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").Value = "'4.238"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").Value = "'0.05136"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.14%  '
$ws.Range("D34").Value = "'0.8056"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +11.34%  '
$ws.Range("D35").Value = "'1.931"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("D36").Value = "'1.167"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("D37").Value = "'2.684"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = "'0.01850"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("D40").Value = '1.175.43'
$ws.Range("E40").Value = '  -6.82%  '
$ws.Range("D41").Value = "'6.203"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("D42").Value = "'0.8933"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.87%  '
$ws.Range("D43").Value = "'72.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'101.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").Value = '2.017.11'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D47").Value = "'0.5179"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").Value = "'1.787"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = "'9.242"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").Value = "'0.9990"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.12%  '
